# Error Calculations and Plots
#
# The source "missing data" sheet simulates random data removal for an
# imputation study. This edit:
#   1. Removes two rows entirely (their IDs "RM 232" and "SC 92" no longer
#      appear), shifting every row below them up.
#   2. Re-randomizes which cells in columns B:F are "missing" (blanked out)
#      vs. populated, for the remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two rows that disappear from the data set ---
# Row 26 ("RM 232") is removed first; everything below shifts up by one,
# so the row that was "SC 92" (originally row 28) is now row 27.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# --- 2. Apply the new "missing data" pattern to the remaining rows ---
# (row numbers below are the FINAL row numbers, after the deletions above)

$ws.Range("E2").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("E5").Value = -5
$ws.Range("D6").Value = -14.2
$ws.Range("E6").Value = -5.7
$ws.Range("F6").Value = 16.43
$ws.Range("D8").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = 16.43
$ws.Range("F11").Value = 17.65
$ws.Range("D12").Value = -14.1
$ws.Range("F12").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("F14").Value = 17.76
$ws.Range("F16").Value = ""
$ws.Range("D17").Value = -14.7
$ws.Range("F17").Value = ""
$ws.Range("D18").Value = -15.2
$ws.Range("D19").Value = ""
$ws.Range("F19").Value = 17.81
$ws.Range("D20").Value = ""
$ws.Range("F21").Value = 16.58
$ws.Range("F22").Value = 16.81
$ws.Range("D23").Value = -13.9
$ws.Range("E24").Value = -8.1
$ws.Range("F25").Value = ""
$ws.Range("F26").Value = ""
$ws.Range("B27").Value = -20.4
$ws.Range("D27").Value = ""
$ws.Range("F27").Value = ""
$ws.Range("E28").Value = ""
$ws.Range("F28").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("E30").Value = -5.7
$ws.Range("F31").Value = 17.18
$ws.Range("B32").Value = ""
